# DE233274: service map upload done
# Rename the single worksheet from the generic default "Sheet1" to a
# descriptive name that matches the template's purpose.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Sales Level Split Upload"
